$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that get newly marked as Milestone "III" complete (E col) with an "X" in F col.
$rows = @(22, 25, 42, 68, 74, 76, 77, 87)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "III"
    $ws.Cells.Item($r, 6).Value = "X"
}

# Re-add the "render to texture" flags on rows 91 and 92 (column E).
$ws.Cells.Item(91, 5).Value = "X"
$ws.Cells.Item(92, 5).Value = "X"

# Update the view/selection state to match the saved workbook.
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F74").Select()
